$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing records ---
# Row 2: password typo fix (Reuts8888! -> Yosi8888!)
$ws.Range("B2").Value = "Yosi8888!"

# Row 3: username fix (yyyyyyy -> tomerne)
$ws.Range("A3").Value = "tomerne"

# Row 4: username fix (yossiyo2 -> omero)
$ws.Range("A4").Value = "omero"

# Row 5: username casing fix (Neriala -> neriala)
$ws.Range("A5").Value = "neriala"

# --- Append new user record ---
$ws.Range("A6").Value = "yossiso1"
$ws.Range("B6").Value = "Tomer12!"
$ws.Range("C6").Value = "207338351"

# Update the active selection to reflect where the user left off editing
$ws.Range("B8").Select() | Out-Null
